$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2447
$ws.Range("I40").Value = 2361.1667
$ws.Range("J40").Value = 2511.375
$ws.Range("K40").Value = 2361.1667
$ws.Range("L40").Value = 2511.375
$ws.Range("M40").Value = -2186.1667
$ws.Range("N40").Value = -2861.375
# Row 64
$ws.Range("H64").Value = 5718.6
$ws.Range("I64").Value = 4532.3335
$ws.Range("J64").Value = 7498
$ws.Range("K64").Value = 4532.3335
$ws.Range("L64").Value = 7498
$ws.Range("M64").Value = -4284.3335
$ws.Range("N64").Value = -7994
# Row 67
$ws.Range("H67").Value = 5718.6
$ws.Range("I67").Value = 4532.3335
$ws.Range("J67").Value = 7498
$ws.Range("K67").Value = 4532.3335
$ws.Range("L67").Value = 7498
$ws.Range("M67").Value = -3674.3335
$ws.Range("N67").Value = -9214
# Row 115
$ws.Range("H115").Value = 6249.5
$ws.Range("I115").Value = 6249.5
$ws.Range("K115").Value = 18748.5
$ws.Range("M115").Value = -17181.5
# Row 116
$ws.Range("H116").Value = 7864.6665
$ws.Range("J116").Value = 10095
$ws.Range("L116").Value = 10095
$ws.Range("N116").Value = -16979
# Row 141
$ws.Range("H141").Value = 18017.8
$ws.Range("J141").Value = 18030.666
$ws.Range("L141").Value = 54091.99800000001
$ws.Range("N141").Value = -64451.99800000001

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 29
$ws.Range("H29").Value = 12000
$ws.Range("J29").Value = 12000
$ws.Range("L29").Value = 12000
$ws.Range("N29").Value = -12616
# Row 45
$ws.Range("H45").Value = 1012
$ws.Range("I45").Value = 1012
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 1012
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -635
$ws.Range("N45").ClearContents()
# Row 61
$ws.Range("H61").Value = 3942.818
$ws.Range("I61").Value = 3942.818
$ws.Range("K61").Value = 3942.818
$ws.Range("M61").Value = -3730.818
# Row 63
$ws.Range("H63").Value = 3487.611
$ws.Range("I63").Value = 1951.9333
$ws.Range("J63").Value = 11166
$ws.Range("K63").Value = 1951.9333
$ws.Range("L63").Value = 11166
$ws.Range("M63").Value = -1265.9333
$ws.Range("N63").Value = -12538
# Row 66
$ws.Range("H66").Value = 3487.611
$ws.Range("I66").Value = 1951.9333
$ws.Range("J66").Value = 11166
$ws.Range("K66").Value = 9759.666499999999
$ws.Range("L66").Value = 55830
$ws.Range("M66").Value = -6327.666499999999
$ws.Range("N66").Value = -62694
# Row 74
$ws.Range("H74").Value = 1927.9
$ws.Range("I74").Value = 1598.1875
$ws.Range("K74").Value = 1598.1875
$ws.Range("M74").Value = -724.1875
# Row 77
$ws.Range("H77").Value = 1927.9
$ws.Range("I77").Value = 1598.1875
$ws.Range("K77").Value = 7990.9375
$ws.Range("M77").Value = -3622.9375
# Row 132
$ws.Range("H132").Value = 1759.8966
$ws.Range("I132").Value = 1749.7142
$ws.Range("K132").Value = 5249.142599999999
$ws.Range("M132").Value = -2719.142599999999
# Row 135
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
# Row 136
$ws.Range("H136").Value = 3942.818
$ws.Range("I136").Value = 3942.818
$ws.Range("K136").Value = 11828.454
$ws.Range("M136").Value = -9278.454000000002

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1170.6666
$ws.Range("I20").Value = 1006
$ws.Range("K20").Value = 1006
$ws.Range("M20").Value = -759
# Row 54
$ws.Range("H54").Value = 7000
$ws.Range("I54").Value = 7000
$ws.Range("K54").Value = 7000
$ws.Range("M54").Value = -6516
# Row 134
$ws.Range("H134").Value = 2499.6
$ws.Range("I134").Value = 2625
$ws.Range("K134").Value = 7875
$ws.Range("M134").Value = -5340

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4062.375
$ws.Range("I31").Value = 3153.923
$ws.Range("K31").Value = 3153.923
$ws.Range("M31").Value = -2858.923
# Row 34
$ws.Range("H34").Value = 4062.375
$ws.Range("I34").Value = 3153.923
$ws.Range("K34").Value = 3153.923
$ws.Range("M34").Value = -2951.923
# Row 51
$ws.Range("H51").Value = 40099
$ws.Range("J51").Value = 40099
$ws.Range("L51").Value = 40099
$ws.Range("N51").Value = -41571
# Row 61
$ws.Range("H61").Value = 40099
$ws.Range("J61").Value = 40099
$ws.Range("L61").Value = 40099
$ws.Range("N61").Value = -40795
# Row 132
$ws.Range("H132").Value = 2283.4285
$ws.Range("I132").Value = 2220.0625
$ws.Range("J132").Value = 2486.2
$ws.Range("K132").Value = 6660.1875
$ws.Range("L132").Value = 7458.599999999999
$ws.Range("M132").Value = -4130.1875
$ws.Range("N132").Value = -12518.6
# Row 134
$ws.Range("H134").Value = 3314.6428
$ws.Range("I134").Value = 3314.6428
$ws.Range("K134").Value = 9943.928400000001
$ws.Range("M134").Value = -7408.928400000001

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 34
$ws.Range("H34").Value = 4403.385
$ws.Range("J34").Value = 6177.222
$ws.Range("L34").Value = 18531.666
$ws.Range("N34").Value = -18699.666

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 12900
$ws.Range("J57").Value = 23000
$ws.Range("L57").Value = 23000
$ws.Range("N57").Value = -24640
# Row 126
$ws.Range("H126").Value = 799.5
$ws.Range("I126").Value = 799.5
$ws.Range("K126").Value = 2398.5
$ws.Range("M126").Value = 71.5
# Row 132
$ws.Range("H132").Value = 3107.5881
$ws.Range("I132").Value = 2925.4614
$ws.Range("J132").Value = 3699.5
$ws.Range("K132").Value = 8776.3842
$ws.Range("L132").Value = 11098.5
$ws.Range("M132").Value = -6246.3842
$ws.Range("N132").Value = -16158.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 6321.6665
$ws.Range("I68").Value = 6321.6665
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 6321.6665
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -5572.6665
$ws.Range("N68").ClearContents()
# Row 71
$ws.Range("H71").Value = 6321.6665
$ws.Range("I71").Value = 6321.6665
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 31608.3325
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -27864.3325
$ws.Range("N71").ClearContents()
# Row 132
$ws.Range("H132").Value = 3669.6667
$ws.Range("I132").Value = 3669.6667
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 11009.0001
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -8479.000100000001
$ws.Range("N132").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 54
$ws.Range("H54").Value = 9000
$ws.Range("I54").Value = 9000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 9000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -8480
$ws.Range("N54").ClearContents()
# Row 95
$ws.Range("H95").Value = 10895.667
$ws.Range("J95").Value = 10895.667
$ws.Range("L95").Value = 10895.667
$ws.Range("N95").Value = -16387.667
# Row 126
$ws.Range("H126").Value = 2271.1428
$ws.Range("I126").Value = 2316.3333
$ws.Range("K126").Value = 6948.999899999999
$ws.Range("M126").Value = -4478.999899999999
# Row 132
$ws.Range("H132").Value = 2425.389
$ws.Range("I132").Value = 2477.1333
$ws.Range("K132").Value = 7431.3999
$ws.Range("M132").Value = -4901.3999
